$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CourtReports")

# Insert a new column before column D (shifts D..W to E..X).
$ws.Range("D1").EntireColumn.Insert()

# The new column takes on column C's width (mirrors the source workbook,
# where the new COURT_REPORT_TAB column is the same width as POM_ITERATION).
$ws.Columns("D").ColumnWidth = $ws.Columns("C").ColumnWidth()

# Populate the new column D header and row2 value
$ws.Range("D1").Value = "COURT_REPORT_TAB"
$ws.Range("D2").Value = "Click"
